# Auto update Excel log
# Appends new sensor/log rows to the "Proximity" and "Camera" sheets.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$row,
        [string]$date,
        [string]$timestamp,
        [string]$hour,
        [string]$location,
        [string]$value,
        [string]$status
    )

    # Reference the previous row so the newly written cells inherit the
    # same (unstyled) look as the rest of the log instead of picking up
    # an automatic "Text" number format style.
    $refRow = $row - 1

    $colA = $ws.Cells.Item($row, 1)
    $colA.NumberFormat = "@"
    $colA.Value = $date
    $colA.Style = $ws.Cells.Item($refRow, 1).Style

    $colB = $ws.Cells.Item($row, 2)
    $colB.NumberFormat = "@"
    $colB.Value = $timestamp
    $colB.Style = $ws.Cells.Item($refRow, 2).Style

    $colC = $ws.Cells.Item($row, 3)
    $colC.NumberFormat = "@"
    $colC.Value = $hour
    $colC.Style = $ws.Cells.Item($refRow, 3).Style

    $ws.Cells.Item($row, 4).Value = $location
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 6).Value = $status
}

# ---- Proximity sheet: append rows 52-54 ----
$wsProximity = $wb.Worksheets.Item("Proximity")

Add-LogRow $wsProximity 52 "2026-02-01" "14:46:52" "14:00" "Living Room Main Door" "EXIT" "User EXITED Living Room Main Door"
Add-LogRow $wsProximity 53 "2026-02-01" "14:46:57" "14:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $wsProximity 54 "2026-02-01" "14:46:58" "14:00" "Living Room Main Door" "EXIT" "User EXITED Living Room Main Door"

# ---- Camera sheet: append rows 35-37 ----
$wsCamera = $wb.Worksheets.Item("Camera")

Add-LogRow $wsCamera 35 "2026-02-01" "14:46:52" "14:00" "Living Room Main Door" "Image Captured" "Active"
Add-LogRow $wsCamera 36 "2026-02-01" "14:46:54" "14:00" "Living Room Main Door" "Image Captured" "Active"
Add-LogRow $wsCamera 37 "2026-02-01" "14:46:57" "14:00" "Living Room Main Door" "Image Received" "Active"
